$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 73: new timesheet entry dated 2024-09-09 (serial 45544), 2.5 hours.
# Copy formatting (date number format) from the cell above (A72) so the new
# date cell picks up the existing style instead of creating a new one.
$ws.Range("A72").Copy()
$ws.Range("A73").PasteSpecial(-4122) | Out-Null
$ws.Range("A73").Value = 45544
$ws.Range("B73").Value = 2.5

# Running total continues the same pattern as the rows above it.
$ws.Range("C73").Formula = "=C72+B73"

# Match the author's final selection/cursor position on the newly added row.
$ws.Range("C73").Select() | Out-Null
